$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = "id"
$ws1.Range("B1").Value = "nama"
$ws1.Range("C1").Value = "item"
$ws1.Range("D1").Value = "kaka"
$ws1.Range("E1").Value = "item"
$ws1.Range("F1").Value = "kaka"
$ws1.Range("G1").Value = "kaka"
[void]$ws1.Range("G8").Select()

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = "id"
[void]$ws2.Range("Q9").Select()
[void]$ws2.Activate()
